# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F
$updates = @{
    3  = 3093
    5  = 2567
    9  = 1315
    13 = 1162
    14 = 337
    16 = 29
    21 = 2356
    22 = 21
    23 = 272
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
